# "Update of the final idea" - rewrite the opening of the Ideation
# paragraph so it reads:
#   "Taking inspiration from both Pong and Pinball, the idea is to
#    create a pinball like game where ..."
# The original sentence started with "The idea is to create ...".
# We keep all of the existing wording (just losing the capital "T" of
# "The", replaced by a lower-case "t" as part of the new lead-in) and
# split the paragraph's single run into three runs, matching the
# target OOXML:
#   1) "Taking inspiration from both Pong and Pinball,"
#   2) " t"
#   3) "he idea is to create a pinball like game ... towards it."

$d = $word.ActiveDocument

# Find the start of the sentence without relying on hard-coded offsets.
$target = $d.Content
$found = $target.Find.Execute("The idea is to create a pinball like game", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $target.Find.Found) {
    throw "Could not find target sentence 'The idea is to create a pinball like game'"
}
$start = $target.Start

# Drop the leading capital "T" of "The" - it gets folded into the new
# " t" run instead.
$rT = $d.Range($start, $start + 1)
$rT.Text = ""

# Insert the new lead-in sentence fragment plus the lower-case "t"
# right before the remaining "he idea is to create ..." text.
$prefix = "Taking inspiration from both Pong and Pinball,"
$joiner = " t"
$rIns = $d.Range($start, $start)
$rIns.InsertBefore($prefix + $joiner)

$lenPrefix = $prefix.Length
$lenJoiner = $joiner.Length

# The text above was inserted as a single run. Force a run boundary
# between the "prefix" piece and the " t" piece by toggling character
# formatting on the " t" slice (on, then back off) - this leaves the
# formatting unchanged but keeps the run split on save, matching the
# three separate <w:r> elements in the target document.
$rJoiner = $d.Range($start + $lenPrefix, $start + $lenPrefix + $lenJoiner)
$rJoiner.Font.Bold = $true
$rJoiner.Font.Bold = $false
